$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: ECs / Efna1 / Epha5 / MuSCs
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.919643
$ws.Range("H2").Value = 53.75892899999999
$ws.Range("I2").Value = 0.8982899767221961
$ws.Range("J2").Value = 0.8982899767221962
$ws.Range("M2").Value = 0.0002903333333333334
$ws.Range("N2").Value = 0.000871
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.005202669684333333
$ws.Range("R2").Value = 0.046824027159
$ws.Range("S2").Value = 0.8982899767221961
$ws.Range("T2").Value = 0.8982899767221962

# Update row 3: FAPs / Efna1 / Epha5 / MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.359006333333333
$ws.Range("H3").Value = 4.077019
$ws.Range("I3").Value = 0.06812533974785755
$ws.Range("J3").Value = 0.06812533974785755
$ws.Range("M3").Value = 0.0002903333333333334
$ws.Range("N3").Value = 0.000871
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.0003945648387777778
$ws.Range("R3").Value = 0.003551083549
$ws.Range("S3").Value = 0.06812533974785755
$ws.Range("T3").Value = 0.06812533974785755

# Update row 4: MuSCs / Efna1 / Epha5 / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.669968
$ws.Range("H4").Value = 2.009904
$ws.Range("I4").Value = 0.03358468352994624
$ws.Range("J4").Value = 0.03358468352994624
$ws.Range("M4").Value = 0.0002903333333333334
$ws.Range("N4").Value = 0.000871
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.0001945140426666667
$ws.Range("R4").Value = 0.001750626384
$ws.Range("S4").Value = 0.03358468352994624
$ws.Range("T4").Value = 0.03358468352994624

# Remove rows 5, 6, 7 (old data no longer present)
$ws.Range("A5:T7").Delete()
